# Auto-generated: apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.671.14'
$ws.Range('E2').Value = '  +0.29%  '

$ws.Range('D3').Value = '1.697.76'
$ws.Range('E3').Value = '  +0.20%  '

$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '''315.85'
$ws.Range('E5').Value = '  -0.15%  '

$ws.Range('D6').Value = '''1.002'
$ws.Range('E6').Value = '  +0.07%  '

$ws.Range('D7').Value = '''0.3927'
$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').Value = '''0.4036'
$ws.Range('E8').Value = '  +0.52%  '

$ws.Range('E9').Value = '  -0.84%  '

$ws.Range('D10').Value = '''1.002'
$ws.Range('E10').Value = '  +0.13%  '

$ws.Range('D11').Value = '''52.97'
$ws.Range('E11').Value = '  -1.86%  '

$ws.Range('D12').Value = '''0.08832'
$ws.Range('E12').Value = '  +0.92%  '

$ws.Range('D13').Value = '''7.462'
$ws.Range('E13').Value = '  +3.60%  '

$ws.Range('D14').Value = '''23.58'
$ws.Range('E14').Value = '  +1.62%  '

$ws.Range('D15').Value = '''8.211'
$ws.Range('E15').Value = '  +7.94%  '

$ws.Range('E16').Value = '  -0.09%  '

$ws.Range('D17').Value = '1.702.44'
$ws.Range('E17').Value = '  +0.29%  '

$ws.Range('D18').Value = '''99.53'
$ws.Range('E18').Value = '  -0.77%  '

$ws.Range('D19').Value = '''0.07027'
$ws.Range('E19').Value = '  -0.24%  '

$ws.Range('D20').Value = '''19.69'
$ws.Range('E20').Value = '  +0.34%  '

$ws.Range('D21').Value = '''7.099'
$ws.Range('E21').Value = '  +3.58%  '

$ws.Range('E22').Value = '  +0.44%  '

$ws.Range('D23').Value = '''14.61'

$ws.Range('D24').Value = '24.652.33'
$ws.Range('E24').Value = '  +0.29%  '

$ws.Range('D25').Value = '''3.114'
$ws.Range('E25').Value = '  +3.39%  '

$ws.Range('D26').Value = '''2.359'
$ws.Range('E26').Value = '  +2.14%  '

$ws.Range('E27').Value = '  +1.18%  '

$ws.Range('D28').Value = '''162.88'
$ws.Range('E28').Value = '  +2.36%  '

$ws.Range('D29').Value = '''8.663'
$ws.Range('E29').Value = '  +15.63%  '

$ws.Range('D30').Value = '''135.50'
$ws.Range('E30').Value = '  +1.45%  '

$ws.Range('D31').Value = '''5.154'
$ws.Range('E31').Value = '  -1.06%  '

$ws.Range('D32').Value = '''0.08958'
$ws.Range('E32').Value = '  +5.03%  '

$ws.Range('D33').Value = '''7.568'
$ws.Range('E33').Value = '  +3.84%  '

$ws.Range('E34').Value = '  -2.69%  '

$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '''1.966'
$ws.Range('E35').Value = '  -0.19%  '

$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '''11.07'
$ws.Range('E36').Value = '  -2.37%  '

$ws.Range('D37').Value = '''0.2755'
$ws.Range('E37').Value = '  +1.24%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.02835'
$ws.Range('E38').Value = '  +3.00%  '

$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '''14.40'
$ws.Range('E39').Value = '  -0.93%  '

$ws.Range('D40').Value = '''0.09117'
$ws.Range('E40').Value = '  +1.08%  '

$ws.Range('D41').Value = '''1.454'
$ws.Range('E41').Value = '  -0.85%  '

$ws.Range('E42').Value = '  -0.26%  '

$ws.Range('D43').Value = '''15.82'
$ws.Range('E43').Value = '  +3.14%  '

$ws.Range('D44').Value = '''0.7159'
$ws.Range('E44').Value = '  -0.27%  '

$ws.Range('D45').Value = '''2.554'
$ws.Range('E45').Value = '  +2.01%  '

$ws.Range('D46').Value = '''4.209'
$ws.Range('E46').Value = '  +0.03%  '

$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('D48').Value = '''1.346'
$ws.Range('E48').Value = '  -0.22%  '

$ws.Range('D49').Value = '''139.77'
$ws.Range('E49').Value = '  -0.94%  '

$ws.Range('D50').Value = '''0.07970'
$ws.Range('E50').Value = '  -0.62%  '

$ws.Range('D51').Value = '''90.41'
$ws.Range('E51').Value = '  +2.50%  '
